$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old row 2 (unit sub-header: Hiver/Eté/Année under (MW)/(GWh)) is removed;
# the remaining data rows shift up by one.
$ws.Rows.Item(2).Delete()

# Row 1 becomes the single header row. Columns A:E are brand-new index/meta
# columns that never had any value or formatting before - make sure they end
# up with the default (unformatted) style.
$ws.Range("A1:E1").ClearFormats()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Columns F:K keep the existing measurement headers, just renamed.
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give the renamed measurement headers (F1:K1) the same font as the rest of
# the header row / table (Arial 9) instead of the plain default font.
$ws.Range("F1:K1").Font.Size = 9

$ws.Range("A2:K2").Select() | Out-Null
